# Update the "Förändrad" (Changed) date column (C) for all data rows.
# The workbook's automatic daily update bumps the serial date value
# in column C from 46075 to 46076 (i.e. 2026-02-22 -> 2026-02-23)
# for every populated row (rows 2 through 33).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 46075) {
        $cell.Value2 = 46076
    }
}
